# Refresh the cryptos list values per the Jan 14 2024 GitHub Actions run.
# (Coin, Link, Price, Volume(1h) columns are B:E; A is just a row index.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.895.76"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").Value = "2.537.87"
$ws.Range("E3").Value = "  -0.83%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.46"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +1.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.35"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +7.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.578"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  +0.67%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.547"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.34"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +2.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0822"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +1.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.64"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -0.91%  "
$ws.Range("E13").Value = "  -0.72%  "
$ws.Range("D14").Value = "2.931.98"
$ws.Range("E14").Value = "  -0.57%  "
$ws.Range("D15").Value = "2.575.80"
$ws.Range("E15").Value = "  +2.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.28"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +7.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.874"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  -0.91%  "
$ws.Range("D18").Value = "42.928.84"
$ws.Range("E18").Value = "  -0.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.17"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +3.38%  "
$ws.Range("E20").Value = "  +0.51%  "
$ws.Range("E21").Value = "  -0.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.71"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -0.78%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "254.30"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -0.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.95"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +0.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.07"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -3.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.73"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -4.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("E28").Value = "  +9.31%  "
$ws.Range("E29").Value = "  -0.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.80"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +4.44%  "
$ws.Range("E31").Value = "  +0.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "157.76"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +3.10%  "
$ws.Range("E33").Value = "  -0.63%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0803"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  +0.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.99"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +7.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.30"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -3.05%  "
$ws.Range("E37").Value = "  -4.25%  "
$ws.Range("E38").Value = "  +0.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "24.33"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +6.29%  "
$ws.Range("E40").Value = "  +0.33%  "
$ws.Range("E41").Value = "  -4.74%  "
$ws.Range("E42").Value = "  +0.63%  "
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0305"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -2.21%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.076.44"
$ws.Range("E45").Value = "  -1.60%  "
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.32"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +1.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.03"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -1.25%  "
$ws.Range("D49").Value = "2.788.24"
$ws.Range("E49").Value = "  -0.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.86"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +2.93%  "
$ws.Range("E51").Value = "  +1.24%  "
